# fujimi-gamebook/checklist.xlsx
# "re-arrange and add missing images"
#
# The Blood Sword block (rows 18-21 before the edit) and the Moonshae Saga
# block (rows 22-27 before the edit, originally at rows 18-23) swap places:
# Blood Sword now comes first (rows 18-21), Moonshae Saga follows
# (rows 22-27). Along the way the wrong "moonshae4.jpeg" image filename is
# corrected to "moonshae4.jpg".
#
# Rows 28-38 keep their data/order - nothing to do there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target contents (after the rearrange) for rows 18-27, columns A-G:
# year, japanese, english, publisher, image, product_type, code

$data = @(
    @(18, 1988, "ブラッド・ソード〈シナリオ1〉勝利の紋章を奪え!", "The Battlepits of Krath", "Fujimi Shobo", "blood_sword1.jpg", "gamebook", "2-1"),
    @(19, 1988, "ブラッド・ソード〈シナリオ2〉魔術王をたおせ!", "The Kingdom of Wyrd", "Fujimi Shobo", "blood_sword2.jpg", "gamebook", "2-2"),
    @(20, 1989, "ブラッド・ソード〈シナリオ3〉悪魔の爪を折れ!", "The Demon's Claw", "Fujimi Shobo", "blood_sword3.jpg", "gamebook", "2-3"),
    @(21, 1989, "ブラッド・ソード〈シナリオ4〉死者の国から還れ!", "Doomwalk", "Fujimi Shobo", "blood_sword4.jpg", "gamebook", "2-4"),
    @(22, 1989, "ムーンシェイ・サーガ〈1〉魔獣よみがえる", "Moonshae Saga 1: Revive the Beast", "Fujimi Shobo", "moonshae1.jpg", "gamebook", "3-1"),
    @(23, 1989, "ムーンシェイ・サーガ〈2〉竪琴と一角獣", "Moonshae Saga 2: Lyre and Unicorn", "Fujimi Shobo", "moonshae2.jpg", "gamebook", "3-2"),
    @(24, 1989, "ムーンシェイ・サーガ〈3〉七人の黒魔術師", "Moonshae Saga 3: Seven Black Magicians", "Fujimi Shobo", "moonshae3.jpg", "gamebook", "3-3"),
    @(25, 1989, "ムーンシェイ・サーガ〈4〉死せる王妃の預言", "Moonshae Saga 4: The Prophecy of the Dead Queen", "Fujimi Shobo", "moonshae4.jpg", "gamebook", "3-4"),
    @(26, 1990, "ムーンシェイ・サーガ〈5〉猫の爪・豹の牙", "Moonshae Saga 5: Cat's Claws, Panther's Fangs", "Fujimi Shobo", "moonshae5.jpg", "gamebook", "3-5"),
    @(27, 1990, "ムーンシェイ・サーガ〈6〉暗黒の解放", "Moonshae Saga 6: The Release of Darkness", "Fujimi Shobo", "moonshae6.jpg", "gamebook", "3-6")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# Widen the "image" column (E) so the longer filenames are readable.
$ws.Columns.Item(5).ColumnWidth = 40

# Move the active selection to E26, matching the author's cursor position.
$ws.Range("E26").Select()
